$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 516.0968
$ws.Range("I33").Value = 433.1111
$ws.Range("K33").Value = 433.1111
$ws.Range("M33").Value = -204.1111
$ws.Range("H70").Value = 1873.8438
$ws.Range("I70").Value = 960
$ws.Range("J70").Value = 3048.7856
$ws.Range("K70").Value = 2880
$ws.Range("L70").Value = 9146.356800000001
$ws.Range("M70").Value = -2610
$ws.Range("N70").Value = -9686.356800000001
$ws.Range("H73").Value = 1873.8438
$ws.Range("I73").Value = 960
$ws.Range("J73").Value = 3048.7856
$ws.Range("K73").Value = 2880
$ws.Range("L73").Value = 9146.356800000001
$ws.Range("M73").Value = -1944
$ws.Range("N73").Value = -11018.3568
$ws.Range("H88").Value = 5972.846
$ws.Range("I88").Value = 6132.727
$ws.Range("J88").Value = 5855.6
$ws.Range("K88").Value = 6132.727
$ws.Range("L88").Value = 5855.6
$ws.Range("M88").Value = -5726.727
$ws.Range("N88").Value = -6667.6
$ws.Range("H91").Value = 5972.846
$ws.Range("I91").Value = 6132.727
$ws.Range("J91").Value = 5855.6
$ws.Range("K91").Value = 6132.727
$ws.Range("L91").Value = 5855.6
$ws.Range("M91").Value = -4728.727
$ws.Range("N91").Value = -8663.6
$ws.Range("H106").Value = 3166.6667
$ws.Range("I106").Value = 500
$ws.Range("J106").Value = 3700
$ws.Range("K106").Value = 500
$ws.Range("L106").Value = 3700
$ws.Range("M106").Value = 131
$ws.Range("N106").Value = -4962
$ws.Range("H113").Value = 4372.5
$ws.Range("I113").Value = 4565.5
$ws.Range("J113").Value = 3890
$ws.Range("K113").Value = 4565.5
$ws.Range("L113").Value = 3890
$ws.Range("M113").Value = -1311.5
$ws.Range("N113").Value = -10398
$ws.Range("H115").Value = 1520
$ws.Range("I115").Value = 533.3333
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 1599.9999
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = -32.99990000000003
$ws.Range("N115").Value = -12134
$ws.Range("H129").Value = 5435901.5
$ws.Range("J129").Value = 1089.7317
$ws.Range("L129").Value = 3269.1951
$ws.Range("N129").Value = -13269.1951
$ws.Range("H137").Value = 2328918.2
$ws.Range("I137").Value = 3452019
$ws.Range("J137").Value = 2495.7856
$ws.Range("K137").Value = 10356057
$ws.Range("L137").Value = 7487.3568
$ws.Range("M137").Value = -10353507
$ws.Range("N137").Value = -12587.3568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8390.683000000001
$ws.Range("I32").Value = 7302
$ws.Range("J32").Value = 24448.75
$ws.Range("K32").Value = 7302
$ws.Range("L32").Value = 24448.75
$ws.Range("M32").Value = -7015
$ws.Range("N32").Value = -25022.75
$ws.Range("H63").Value = 2223.111
$ws.Range("H66").Value = 2223.111
$ws.Range("H97").Value = 789.125
$ws.Range("I97").Value = 761.7646999999999
$ws.Range("K97").Value = 761.7646999999999
$ws.Range("M97").Value = -265.7646999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3127.9048
$ws.Range("I99").Value = 1933.7333
$ws.Range("J99").Value = 6113.3335
$ws.Range("K99").Value = 1933.7333
$ws.Range("L99").Value = 6113.3335
$ws.Range("M99").Value = -435.7333000000001
$ws.Range("N99").Value = -9109.333500000001
$ws.Range("H107").Value = 1948.8572
$ws.Range("I107").Value = 579.1818
$ws.Range("J107").Value = 6971
$ws.Range("K107").Value = 579.1818
$ws.Range("L107").Value = 6971
$ws.Range("M107").Value = 1340.8182
$ws.Range("N107").Value = -10811
$ws.Range("H126").Value = 3127.9048
$ws.Range("I126").Value = 1933.7333
$ws.Range("J126").Value = 6113.3335
$ws.Range("K126").Value = 5801.199900000001
$ws.Range("L126").Value = 18340.0005
$ws.Range("M126").Value = -3331.199900000001
$ws.Range("N126").Value = -23280.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 76.064514
$ws.Range("I2").Value = 90.22727
$ws.Range("J2").Value = 41.444443
$ws.Range("K2").Value = 90.22727
$ws.Range("L2").Value = 41.444443
$ws.Range("M2").Value = 22.77273
$ws.Range("N2").Value = -267.444443
$ws.Range("H80").Value = 3735.4285
$ws.Range("I80").Value = 3168.3333
$ws.Range("J80").Value = 4160.75
$ws.Range("K80").Value = 3168.3333
$ws.Range("L80").Value = 4160.75
$ws.Range("M80").Value = -2170.3333
$ws.Range("N80").Value = -6156.75
$ws.Range("H83").Value = 3735.4285
$ws.Range("I83").Value = 3168.3333
$ws.Range("J83").Value = 4160.75
$ws.Range("K83").Value = 15841.6665
$ws.Range("L83").Value = 20803.75
$ws.Range("M83").Value = -10849.6665
$ws.Range("N83").Value = -30787.75
$ws.Range("H102").Value = 34212.53
$ws.Range("I102").Value = 2350.1
$ws.Range("J102").Value = 87316.586
$ws.Range("K102").Value = 2350.1
$ws.Range("L102").Value = 87316.586
$ws.Range("M102").Value = -728.0999999999999
$ws.Range("N102").Value = -90560.586

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 111112670
$ws.Range("I22").Value = 200000400
$ws.Range("J22").Value = 3012.5
$ws.Range("K22").Value = 200000400
$ws.Range("L22").Value = 3012.5
$ws.Range("M22").Value = -200000105
$ws.Range("N22").Value = -3602.5
$ws.Range("H27").Value = 111112670
$ws.Range("I27").Value = 200000400
$ws.Range("J27").Value = 3012.5
$ws.Range("K27").Value = 200000400
$ws.Range("L27").Value = 3012.5
$ws.Range("M27").Value = -200000293
$ws.Range("N27").Value = -3226.5
$ws.Range("H61").Value = 2739.8
$ws.Range("I61").Value = 933
$ws.Range("J61").Value = 5450
$ws.Range("K61").Value = 933
$ws.Range("L61").Value = 5450
$ws.Range("M61").Value = -731
$ws.Range("N61").Value = -5854
$ws.Range("H93").Value = 1668.7693
$ws.Range("I93").Value = 965.1111
$ws.Range("J93").Value = 3252
$ws.Range("K93").Value = 965.1111
$ws.Range("L93").Value = 3252
$ws.Range("M93").Value = 282.8889
$ws.Range("N93").Value = -5748
$ws.Range("H100").Value = 2837.7144
$ws.Range("J100").Value = 2837.7144
$ws.Range("L100").Value = 2837.7144
$ws.Range("N100").Value = -3919.7144
$ws.Range("H113").Value = 2739.8
$ws.Range("I113").Value = 933
$ws.Range("J113").Value = 5450
$ws.Range("K113").Value = 933
$ws.Range("L113").Value = 5450
$ws.Range("M113").Value = 1237
$ws.Range("N113").Value = -9790

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1424.45
$ws.Range("I81").Value = 1179.5
$ws.Range("J81").Value = 1669.4
$ws.Range("K81").Value = 2359
$ws.Range("L81").Value = 3338.8
$ws.Range("M81").Value = -1298
$ws.Range("N81").Value = -5460.8
$ws.Range("H84").Value = 1424.45
$ws.Range("I84").Value = 1179.5
$ws.Range("J84").Value = 1669.4
$ws.Range("K84").Value = 11795
$ws.Range("L84").Value = 16694
$ws.Range("M84").Value = -6491
$ws.Range("N84").Value = -27302
$ws.Range("H107").Value = 636.03705
$ws.Range("I107").Value = 262.65
$ws.Range("J107").Value = 1702.8572
$ws.Range("K107").Value = 787.9499999999999
$ws.Range("L107").Value = 5108.571599999999
$ws.Range("M107").Value = 1132.05
$ws.Range("N107").Value = -8948.571599999999
$ws.Range("H128").Value = 39900
$ws.Range("J128").Value = 39900
$ws.Range("L128").Value = 39900
$ws.Range("N128").Value = -49860
$ws.Range("H132").Value = 148232.03
$ws.Range("I132").Value = 167346.25
$ws.Range("J132").Value = 4875.375
$ws.Range("K132").Value = 502038.75
$ws.Range("L132").Value = 14626.125
$ws.Range("M132").Value = -499508.75
$ws.Range("N132").Value = -19686.125
